$d = $word.ActiveDocument

# The name heading ("Dheeraj Chand") is the first paragraph in the body.
# We need to add a new, centered paragraph right after it that contains
# the contact information line (this is the "short" resume contact-info
# paragraph that was missing).
$namePara = $d.Paragraphs.First
$nameRange = $namePara.Range

# Create a new paragraph right after the name paragraph. Word clones the
# preceding paragraph/run formatting (bold, 28-half-point size, centered
# alignment) onto this new, still-empty paragraph.
$nameRange.InsertParagraphAfter()
$contactPara = $d.Paragraphs(2)

# Replace the (inherited-formatting) contents of the new paragraph with a
# clean OOXML fragment: centered alignment, but a plain, unformatted run
# for the contact-info text -- matching a freshly authored paragraph
# rather than one that drags along the name's bold/large-font run
# properties.
$contactXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:pPr><w:jc w:val="center"/></w:pPr>' + `
              '<w:r><w:t>202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX</w:t></w:r>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$contactPara.Range.InsertXML($contactXml)

Write-Output "Inserted contact info paragraph after the name heading."
